# Populate the "key"/"value" settings table on both worksheets
# (aurix_app / aurix_fbl) with the login + build metadata rows, then
# leave the active-cell selection where the author left off.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

function Fill-SettingsSheet($ws) {
    # --- column A (keys) for the login row, then column B (values) ---
    $ws.Cells.Item(2, 1).Value = "login.username"
    $ws.Cells.Item(3, 1).Value = "login.password"
    $ws.Cells.Item(2, 2).Value = "zhoutao"
    $ws.Cells.Item(3, 2).Value = "Ztb547611679!!"

    # --- remaining keys in column A ---
    $ws.Cells.Item(4, 1).Value = "version"
    $ws.Cells.Item(6, 1).Value = "semantic_version"
    $ws.Cells.Item(5, 1).Value = "general_setting"
    $ws.Cells.Item(7, 1).Value = "software_part_number"
    $ws.Cells.Item(8, 1).Value = "software_YMP_version"
    $ws.Cells.Item(9, 1).Value = "dependencies"
    $ws.Cells.Item(10, 1).Value = "file_upload_ODX_F"
    $ws.Cells.Item(11, 1).Value = "file_upload_flashware"

    # --- remaining values in column B ---
    $ws.Cells.Item(4, 2).Value = "1.2.0"
    $ws.Cells.Item(5, 2).Value = "Create new"
    $ws.Cells.Item(7, 2).Value = "A1234567890"

    # "01/01/00" must stay literal text, not become a date serial number
    $ws.Cells.Item(8, 2).NumberFormat = "@"
    $ws.Cells.Item(8, 2).Value = "01/01/00"
    $ws.Cells.Item(8, 2).ClearFormats()

    $ws.Cells.Item(9, 2).Value = "A1234567890 01/01/00"
    $ws.Cells.Item(6, 2).Value = "1.0.0-alph"
    $ws.Cells.Item(10, 2).Value = "C:\Users\zhoutao\Downloads\1224 Automatic uploading MBOS test\0009047908_253629_FULL_IDC_GEN6PG_C_STAR35_AURIX_FBL.odx-f"
    $ws.Cells.Item(11, 2).Value = "C:\Users\zhoutao\Downloads\1224 Automatic uploading MBOS test\0019023738_254669.hex"
}

Fill-SettingsSheet $ws1
Fill-SettingsSheet $ws2

# Restore the selections the author left on each sheet
$ws1.Activate()
$ws1.Range("B19").Select()

$ws2.Activate()
$ws2.Range("B22").Select()
